$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Burndown")

# Column C (Base Line) - shift all values up by 5
$ws.Range("C5").Value = 45
$ws.Range("C6").Value = 40
$ws.Range("C7").Value = 35
$ws.Range("C8").Value = 30
$ws.Range("C9").Value = 25
$ws.Range("C10").Value = 20
$ws.Range("C11").Value = 15
$ws.Range("C12").Value = 10
$ws.Range("C13").Value = 5

# Column D (Horas Trabalhadas) - new actuals
$ws.Range("D5").Value = 45
$ws.Range("D6").Value = 38.5
$ws.Range("D7").Value = 34
$ws.Range("D8").Value = 28
$ws.Range("D9").Value = 25
$ws.Range("D10").Value = 20.5
$ws.Range("D11").Value = 15.5
$ws.Range("D12").Value = 10
$ws.Range("D13").Value = 8.5
